$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.554.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.53%  '
$ws.Range("E3").Value = '  +2.95%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '401.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.26%  '
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +5.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0879'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.12%  '
$ws.Range("E13").Value = '  +2.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.99%  '
$ws.Range("E16").Value = '  +8.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.184.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '54.484.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.05%  '
$ws.Range("E20").Value = '  +4.48%  '
$ws.Range("E21").Value = '  +3.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0993'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '274.60'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.03'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.55%  '
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  +4.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0505'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +12.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.89%  '
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '50.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("E37").Value = '  +7.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.82%  '
$ws.Range("E41").Value = '  +3.28%  '
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '130.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.117'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.090.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0346'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0509'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.27%  '
